$p = $ppt.ActivePresentation

function Set-ParagraphText {
    param($Slide, $ShapeIndex, $ParaIndex, $OldText, $NewText)
    $shp = $Slide.Shapes.Item($ShapeIndex)
    $tr = $shp.TextFrame.TextRange
    $para = $tr.Paragraphs($ParaIndex, 1)
    if ($para.Text -ne $OldText) {
        Write-Host "WARNING: unexpected text at shape=" $ShapeIndex "para=" $ParaIndex ":" $para.Text "expected:" $OldText
    }
    $para.Text = $NewText
}

Set-ParagraphText $p.Slides.Item(10) 2 2 "Nº de días de lluvia en el mes/año" "Nº de días de lluvia en el mes"
Set-ParagraphText $p.Slides.Item(10) 3 2 " Nº de días de nieve en el mes/año" " Nº de días de nieve en el mes"
Set-ParagraphText $p.Slides.Item(11) 1 3 "Nº de días de granizo en el mes/año" "Nº de días de granizo en el mes"
Set-ParagraphText $p.Slides.Item(11) 2 2 "Nº de días de tormenta en el mes/año" "Nº de días de tormenta en el mes"
Set-ParagraphText $p.Slides.Item(12) 1 2 "Nº de días de niebla en el mes/año" "Nº de días de niebla en el mes"
Set-ParagraphText $p.Slides.Item(12) 2 2 "Nº de días despejados en el mes/año" "Nº de días despejados en el mes"
Set-ParagraphText $p.Slides.Item(13) 1 2 "Nº de días nubosos en el mes/año" "Nº de días nubosos en el mes"
Set-ParagraphText $p.Slides.Item(13) 2 2 "Nº de días cubiertos en el mes/año" "Nº de días cubiertos en el mes"
Set-ParagraphText $p.Slides.Item(16) 2 2 "Temperatura media mensual/anual (°C)" "Temperatura media mensual (°C)"
Set-ParagraphText $p.Slides.Item(16) 3 2 "Temperatura media mensual/anual de las máximas (°C)" "Temperatura media mensual de las máximas (°C)"
Set-ParagraphText $p.Slides.Item(17) 1 2 "Temperatura media mensual/anual de las mínimas (°C)" "Temperatura media mensual de las mínimas (°C)"
Set-ParagraphText $p.Slides.Item(17) 2 2 "Temperatura máxima absoluta del mes/año y fecha (°C)" "Temperatura máxima absoluta del mes (°C)"
Set-ParagraphText $p.Slides.Item(18) 1 2 "Temperatura mínima absoluta del mes/año y fecha (°C)" "Temperatura mínima absoluta del mes (°C)"
Set-ParagraphText $p.Slides.Item(18) 2 2 "Temperatura mínima más alta del mes/año (°C)" "Temperatura mínima más alta del mes (°C)"
Set-ParagraphText $p.Slides.Item(19) 1 2 "Temperatura máxima más baja del mes/año (°C)" "Temperatura máxima más baja del mes (°C)"
Set-ParagraphText $p.Slides.Item(19) 2 2 "Nº días mes/año de temperatura máxima mayor o igual que 30°C" "Nº días mes de temperatura máxima mayor o igual que 30°C"
Set-ParagraphText $p.Slides.Item(20) 1 2 "Nº días mes/año de temperatura mínima menor o igual que 0°C" "Nº días mes de temperatura mínima menor o igual que 0°C"
Set-ParagraphText $p.Slides.Item(21) 2 2 "Nº días de velocidad del viento mayor o igual a 55 Km/h en el mes/año" "Nº días de velocidad del viento mayor o igual a 55 Km/h en el mes"
Set-ParagraphText $p.Slides.Item(21) 3 2 "Nº días de velocidad del viento mayor o igual a 91 Km/h en el mes/año" "Nº días de velocidad del viento mayor o igual a 91 Km/h en el mes"
Set-ParagraphText $p.Slides.Item(7) 2 2 "Precipitación total mensual/anual (mm)" "Precipitación total mensual (mm)"
Set-ParagraphText $p.Slides.Item(7) 3 2 "Precipitación máxima diaria del mes/año y fecha (mm)" "Precipitación máxima diaria del mes (mm)"
Set-ParagraphText $p.Slides.Item(8) 1 2 "Nº de días de precipitación apreciable (≥ 0,1 mm) en el mes/año" "Nº de días de precipitación apreciable (≥ 0,1 mm) en el mes"
Set-ParagraphText $p.Slides.Item(8) 2 2 "Nº de días de precipitación mayor o igual que 1mm en el mes/año" "Nº de días de precipitación mayor o igual que 1mm en el mes"
Set-ParagraphText $p.Slides.Item(9) 1 2 "Nº de días de precipitación mayor o igual que 10mm en el mes/año" "Nº de días de precipitación mayor o igual que 10mm en el mes"
Set-ParagraphText $p.Slides.Item(9) 2 2 "Nº de días de precipitación mayor o igual que 30mm en el mes/año" "Nº de días de precipitación mayor o igual que 30mm en el mes"
